# Updates Cntn1-Notch1.xlsx LR-pair data with new TPM-based values.
# - Rewrites existing rows 2-6 (FAPs sending cluster) with recomputed stats
# - Appends new rows 7-11 for the MuSCs sending cluster (same ligand/receptor,
#   all 5 target clusters), extending the used range to A1:T11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cntn1"
$ws.Range("C2").Value = "Notch1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07031566666666667
$ws.Range("H2").Value = 0.210947
$ws.Range("I2").Value = 0.3500709860127268
$ws.Range("J2").Value = 0.446885632088942
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 38.10639333333333
$ws.Range("N2").Value = 114.31918
$ws.Range("O2").Value = 0.3831479157160237
$ws.Range("P2").Value = 0.4159903984418967
$ws.Range("Q2").Value = 2.679476451495555
$ws.Range("R2").Value = 24.11528806346
$ws.Range("S2").Value = 0.1341289686434296
$ws.Range("T2").Value = 0.1859001321506378

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cntn1"
$ws.Range("C3").Value = "Notch1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07031566666666667
$ws.Range("H3").Value = 0.210947
$ws.Range("I3").Value = 0.3500709860127268
$ws.Range("J3").Value = 0.446885632088942
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.149483999999999
$ws.Range("N3").Value = 27.448452
$ws.Range("O3").Value = 0.09199521176963764
$ws.Range("P3").Value = 0.09988081163714851
$ws.Range("Q3").Value = 0.6433520671159999
$ws.Range("R3").Value = 5.790168604044
$ws.Range("S3").Value = 0.03220485449264666
$ws.Range("T3").Value = 0.04463529964202367

# Row 4: FAPs -> Inflammatory-Mac
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cntn1"
$ws.Range("C4").Value = "Notch1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07031566666666667
$ws.Range("H4").Value = 0.210947
$ws.Range("I4").Value = 0.3500709860127268
$ws.Range("J4").Value = 0.446885632088942
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 15.023598
$ws.Range("N4").Value = 45.070794
$ws.Range("O4").Value = 0.1510575983904562
$ws.Range("P4").Value = 0.1640058785774412
$ws.Range("Q4").Value = 1.056394309102
$ws.Range("R4").Value = 9.507548781918
$ws.Range("S4").Value = 0.0528808824132615
$ws.Range("T4").Value = 0.07329187071438209

# Row 5: FAPs -> MuSCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cntn1"
$ws.Range("C5").Value = "Notch1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.07031566666666667
$ws.Range("H5").Value = 0.210947
$ws.Range("I5").Value = 0.3500709860127268
$ws.Range("J5").Value = 0.446885632088942
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 23.556204
$ws.Range("N5").Value = 47.112408
$ws.Range("O5").Value = 0.236850294013169
$ws.Range("P5").Value = 0.1714350065796238
$ws.Range("Q5").Value = 1.656370188396
$ws.Range("R5").Value = 9.938221130376
$ws.Range("S5").Value = 0.08291441596259432
$ws.Range("T5").Value = 0.07661184127750713

# Row 6: FAPs -> Resolving-Mac
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cntn1"
$ws.Range("C6").Value = "Notch1"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.07031566666666667
$ws.Range("H6").Value = 0.210947
$ws.Range("I6").Value = 0.3500709860127268
$ws.Range("J6").Value = 0.446885632088942
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 13.62041
$ws.Range("N6").Value = 40.86123000000001
$ws.Range("O6").Value = 0.1369489801107134
$ws.Range("P6").Value = 0.1486879047638899
$ws.Range("Q6").Value = 0.9577282094233334
$ws.Range("R6").Value = 8.619553884810001
$ws.Range("S6").Value = 0.04794186450079475
$ws.Range("T6").Value = 0.06644648830439134

# Row 7: MuSCs -> ECs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Cntn1"
$ws.Range("C7").Value = "Notch1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.1305455
$ws.Range("H7").Value = 0.261091
$ws.Range("I7").Value = 0.6499290139872732
$ws.Range("J7").Value = 0.5531143679110581
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 38.10639333333333
$ws.Range("N7").Value = 114.31918
$ws.Range("O7").Value = 0.3831479157160237
$ws.Range("P7").Value = 0.4159903984418967
$ws.Range("Q7").Value = 4.974618170896666
$ws.Range("R7").Value = 29.84770902538
$ws.Range("S7").Value = 0.2490189470725942
$ws.Range("T7").Value = 0.2300902662912589

# Row 8: MuSCs -> FAPs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Cntn1"
$ws.Range("C8").Value = "Notch1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.1305455
$ws.Range("H8").Value = 0.261091
$ws.Range("I8").Value = 0.6499290139872732
$ws.Range("J8").Value = 0.5531143679110581
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 9.149483999999999
$ws.Range("N8").Value = 27.448452
$ws.Range("O8").Value = 0.09199521176963764
$ws.Range("P8").Value = 0.09988081163714851
$ws.Range("Q8").Value = 1.194423963522
$ws.Range("R8").Value = 7.166543781132001
$ws.Range("S8").Value = 0.05979035727699098
$ws.Range("T8").Value = 0.05524551199512485

# Row 9: MuSCs -> Inflammatory-Mac
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Cntn1"
$ws.Range("C9").Value = "Notch1"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.1305455
$ws.Range("H9").Value = 0.261091
$ws.Range("I9").Value = 0.6499290139872732
$ws.Range("J9").Value = 0.5531143679110581
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 15.023598
$ws.Range("N9").Value = 45.070794
$ws.Range("O9").Value = 0.1510575983904562
$ws.Range("P9").Value = 0.1640058785774412
$ws.Range("Q9").Value = 1.961263112709
$ws.Range("R9").Value = 11.767578676254
$ws.Range("S9").Value = 0.0981767159771947
$ws.Range("T9").Value = 0.09071400786305914

# Row 10: MuSCs -> MuSCs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Cntn1"
$ws.Range("C10").Value = "Notch1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1305455
$ws.Range("H10").Value = 0.261091
$ws.Range("I10").Value = 0.6499290139872732
$ws.Range("J10").Value = 0.5531143679110581
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 23.556204
$ws.Range("N10").Value = 47.112408
$ws.Range("O10").Value = 0.236850294013169
$ws.Range("P10").Value = 0.1714350065796238
$ws.Range("Q10").Value = 3.075156429282
$ws.Range("R10").Value = 12.300625717128
$ws.Range("S10").Value = 0.1539358780505747
$ws.Range("T10").Value = 0.09482316530211672

# Row 11: MuSCs -> Resolving-Mac
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Cntn1"
$ws.Range("C11").Value = "Notch1"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1305455
$ws.Range("H11").Value = 0.261091
$ws.Range("I11").Value = 0.6499290139872732
$ws.Range("J11").Value = 0.5531143679110581
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 13.62041
$ws.Range("N11").Value = 40.86123000000001
$ws.Range("O11").Value = 0.1369489801107134
$ws.Range("P11").Value = 0.1486879047638899
$ws.Range("Q11").Value = 1.778083233655
$ws.Range("R11").Value = 10.66849940193
$ws.Range("S11").Value = 0.08900711560991864
$ws.Range("T11").Value = 0.08224141645949855
